$d = $word.ActiveDocument

$d.Content.Find.Execute("Terrance Archie", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Delvin Braxton", 2)

$d.Content.Find.Execute("careerpath5498@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4039forscythe@gmail.com", 2)

$d.Content.Find.Execute("45-44 42nd St New York NY 11104", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1331 W Lunt Ave Chicago IL 60626", 2)

$d.Content.Find.Execute("Manhattan College", $true, $false, $false, $false, $false,
                         $true, 1, $false, "University of Illinois at Chicago", 2)
